# "Generate Report for Archive"
#
# Localization status report refresh:
#   - The two translation-status cells ("Ready for handoff") on each
#     language sheet / the Overview sheet flip to "In Translation".
#   - The Status-ish column(s) that held that text are narrowed to fit
#     the new (shorter) text.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: columns E (zh-cn) and F (de-de), rows 2-3 ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# Shrink the now-narrower Status columns to match the shorter text.
# (ColumnWidth is snapped to the host's pixel grid the same way real
# Excel snaps displayed column widths, so we feed it the input that
# lands on the grid point closest to the refreshed report's width.)
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn sheet: Status column C, rows 2-3 ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

# --- de-de sheet: Status column C, rows 2-3 ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
